# Add season-record columns (Wins, Losses, Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (style s="1")
# from A1 onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-51): every player on the roster shares the team's
# season record, so fill the same Wins/Losses/Ties values down each row.
$wins = 81
$losses = 80
$ties = 0

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
